# Update existing monthly values for rows 194-199 (columns J, M, Q, W, X)
# and append a new monthly row (row 200) for 01-07-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    194 = 2374
    195 = 2327
    196 = 2279
    197 = 2200
    198 = 2137
    199 = 2074
}

foreach ($row in $updates.Keys) {
    $val = $updates[$row]
    $ws.Range("J$row").Value = $val
    $ws.Range("M$row").Value = $val
    $ws.Range("Q$row").Value = -$val
    $ws.Range("W$row").Value = $val
    $ws.Range("X$row").Value = -$val
}

# Append new row 200 with the 01-07-2021 monthly figures
$newRow = 200
$newVal = 2007

# Build the date label as a formula result elsewhere and paste the raw
# value into A200 so that Excel stores it as plain text (matching the
# existing "dd-mm-yyyy" text labels in column A) instead of auto-
# converting the string into a date serial number.
$ws.Range("Z1").Formula = "=""01-07-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A$newRow").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents()

foreach ($col in @("B","C","D","E","F","G","H","I","K","L","N","O","P","R","S","T","U","V")) {
    $ws.Range("$col$newRow").Value = 0
}

$ws.Range("J$newRow").Value = $newVal
$ws.Range("M$newRow").Value = $newVal
$ws.Range("Q$newRow").Value = -$newVal
$ws.Range("W$newRow").Value = $newVal
$ws.Range("X$newRow").Value = -$newVal
